$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 457, pushing the existing rows 457-516 down to 459-518.
$ws.Rows("457:458").Insert()

# Row 457 (new): Primera, weekly update
$ws.Cells.Item(457,1).Value = 7
$ws.Cells.Item(457,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(457,3).Value = "Ñuble"
$ws.Cells.Item(457,4).Value = 45124
$ws.Cells.Item(457,5).Value = 16
$ws.Cells.Item(457,6).Value = 100112009
$ws.Cells.Item(457,7).Value = "Acelga"
$ws.Cells.Item(457,8).Value = "Sin especificar"
$ws.Cells.Item(457,9).Value = "Primera"
$ws.Cells.Item(457,10).Value = 300
$ws.Cells.Item(457,11).Value = 600
$ws.Cells.Item(457,12).Value = 700
$ws.Cells.Item(457,13).Value = 650
$ws.Cells.Item(457,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(457,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(457,16).Value = 650
$ws.Cells.Item(457,17).Value = 1
$ws.Cells.Item(457,18).Value = "Hortaliza"

# Row 458 (new): Segunda, weekly update
$ws.Cells.Item(458,1).Value = 7
$ws.Cells.Item(458,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(458,3).Value = "Ñuble"
$ws.Cells.Item(458,4).Value = 45124
$ws.Cells.Item(458,5).Value = 16
$ws.Cells.Item(458,6).Value = 100112009
$ws.Cells.Item(458,7).Value = "Acelga"
$ws.Cells.Item(458,8).Value = "Sin especificar"
$ws.Cells.Item(458,9).Value = "Segunda"
$ws.Cells.Item(458,10).Value = 200
$ws.Cells.Item(458,11).Value = 500
$ws.Cells.Item(458,12).Value = 500
$ws.Cells.Item(458,13).Value = 500
$ws.Cells.Item(458,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(458,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(458,16).Value = 500
$ws.Cells.Item(458,17).Value = 1
$ws.Cells.Item(458,18).Value = "Hortaliza"
